$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.813.32'
$ws.Range("E2").Value = '  +4.95%  '
$ws.Range("D3").Value = '2.280.92'
$ws.Range("E3").Value = '  +3.24%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.47'
$ws.Range("E5").Value = '  +1.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.642'
$ws.Range("E6").Value = '  +3.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.87'
$ws.Range("E7").Value = '  +9.09%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.431'
$ws.Range("E9").Value = '  +6.99%  '
$ws.Range("E10").Value = '  +16.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.65'
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.44'
$ws.Range("E12").Value = '  +18.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("E13").Value = '  +0.71%  '
$ws.Range("D14").Value = '2.621.12'
$ws.Range("E14").Value = '  +3.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.94'
$ws.Range("E15").Value = '  +3.16%  '
$ws.Range("E16").Value = '  +5.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.835'
$ws.Range("E17").Value = '  +5.03%  '
$ws.Range("D18").Value = '2.275.34'
$ws.Range("E18").Value = '  +2.97%  '
$ws.Range("D19").Value = '43.716.45'
$ws.Range("E19").Value = '  +4.96%  '
$ws.Range("D20").Value = '0.0₃0994'
$ws.Range("E20").Value = '  +10.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.26'
$ws.Range("E21").Value = '  +2.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.18'
$ws.Range("E22").Value = '  +1.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '261.85'
$ws.Range("E23").Value = '  +7.88%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.52'
$ws.Range("E25").Value = '  +7.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.33'
$ws.Range("E26").Value = '  +2.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.20'
$ws.Range("E27").Value = '  +5.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '173.02'
$ws.Range("E28").Value = '  +2.11%  '
$ws.Range("E29").Value = '  +6.98%  '
$ws.Range("E30").Value = '  -1.89%  '
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.81'
$ws.Range("E32").Value = '  +8.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.125'
$ws.Range("E33").Value = '  +2.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0690'
$ws.Range("E34").Value = '  +6.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.09'
$ws.Range("E35").Value = '  +1.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.79'
$ws.Range("E36").Value = '  +3.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.84'
$ws.Range("E37").Value = '  +8.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.82'
$ws.Range("E38").Value = '  +7.83%  '
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("E40").Value = '  +4.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E42").Value = '  -1.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0983'
$ws.Range("E44").Value = '  +2.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.84'
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("E48").Value = '  +0.94%  '
$ws.Range("D49").Value = '1.480.39'
$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("E50").Value = '  +7.31%  '
$ws.Range("E51").Value = '  -12.42%  '

# Row reshuffle for FTXToken / InjectiveProtocol / Celestia (rows 43, 45, 46)
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.53'
$ws.Range("E43").Value = '  +6.92%  '
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.52'
$ws.Range("E45").Value = '  +22.70%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.49'
$ws.Range("E46").Value = '  +2.34%  '
